$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "5.020", "0.9975") keep their exact text representation,
# matching the inline-string cells in the target workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.333.24"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "1.846.01"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("D4").Value = "0.9975"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").Value = "240.07"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").Value = "0.6261"
$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("D7").Value = "0.9985"
$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("D8").Value = "0.07591"
$ws.Range("E8").Value = "  -1.31%  "

$ws.Range("E9").Value = "  -1.51%  "

$ws.Range("D10").Value = "24.68"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("D11").Value = "0.07733"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").Value = "5.020"
$ws.Range("E12").Value = "  -0.07%  "

$ws.Range("D13").Value = "0.6781"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").Value = "0.00001060"
$ws.Range("E14").Value = "  -2.33%  "

$ws.Range("D15").Value = "82.87"
$ws.Range("E15").Value = "  -1.04%  "

$ws.Range("D16").Value = "6.114"
$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").Value = "29.349.96"
$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").Value = "227.63"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("E19").Value = "  -1.11%  "

$ws.Range("D20").Value = "0.9984"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D21").Value = "7.475"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "0.9984"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").Value = "158.46"
$ws.Range("E23").Value = "  +0.73%  "

$ws.Range("D24").Value = "0.1380"
$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("D25").Value = "8.422"
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("D26").Value = "17.63"
$ws.Range("E26").Value = "  -0.28%  "

$ws.Range("D27").Value = "1.433"
$ws.Range("E27").Value = "  +8.84%  "

$ws.Range("D28").Value = "1.458"
$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("D29").Value = "0.05602"
$ws.Range("E29").Value = "  -2.03%  "

$ws.Range("D30").Value = "4.098"
$ws.Range("E30").Value = "  -0.44%  "

$ws.Range("D31").Value = "4.062"
$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").Value = "1.828"
$ws.Range("E33").Value = "  -1.21%  "

$ws.Range("E34").Value = "  -1.92%  "

$ws.Range("D35").Value = "2.580"
$ws.Range("E35").Value = "  -0.31%  "

$ws.Range("D36").Value = "0.01798"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").Value = "1.226.48"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").Value = "2.721"
$ws.Range("E38").Value = "  -2.06%  "

$ws.Range("D39").Value = "6.351"
$ws.Range("E39").Value = "  -1.51%  "

$ws.Range("D40").Value = "0.8976"
$ws.Range("E40").Value = "  -1.45%  "

$ws.Range("D41").Value = "0.9981"
$ws.Range("E41").Value = "  -0.25%  "

$ws.Range("D42").Value = "101.44"
$ws.Range("E42").Value = "  -0.38%  "

$ws.Range("D43").Value = "65.43"
$ws.Range("E43").Value = "  -1.15%  "

$ws.Range("D44").Value = "0.00000000121"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").Value = "7.205"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("D46").Value = "0.3984"
$ws.Range("E46").Value = "  -1.05%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "1.688"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.002"
$ws.Range("E48").Value = "  -0.44%  "

$ws.Range("E49").Value = "  +1.37%  "

$ws.Range("D50").Value = "0.05695"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("D51").Value = "0.4620"
$ws.Range("E51").Value = "  -0.18%  "

# Reset the style reference on column D back to Normal/General so the
# saved cells do not carry a stray style index (keeps styles.xml clean)
# while the stored value type remains text.
$ws.Range("D2:D51").Style = "Normal"

